# Adds two new match rows (108, 109) to the Paraguay Primera Division 2023
# results sheet, matching the source diff:
#   row 108 -> Indice 107, Guairena 0 x 1 Olimpia Asuncion
#   row 109 -> Indice 108, Cerro Porteno 2 x 0 Sp. Luqueno

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 108 -------------------------------------------------------------

# Column A ("Indice") carries the bordered/bold/centered header-style (s="1")
# used throughout the sheet. Pull the format from the row above via
# PasteSpecial (xlPasteFormats) so no new style entries get minted, then set
# the value.
$ws.Range("A107").Copy() | Out-Null
$ws.Range("A108").PasteSpecial(-4122) | Out-Null
$ws.Range("A108").Value = 107

# Columns B/C/D repeat constant text values ("paraguay" / "primera-division"
# / "2023") on every data row. Copy them wholesale (value + format) from the
# row above instead of re-typing the strings, which keeps "2023" as text
# instead of Excel auto-coercing it to a number.
$ws.Range("B107").Copy() | Out-Null
$ws.Range("B108").PasteSpecial(-4104) | Out-Null
$ws.Range("C107").Copy() | Out-Null
$ws.Range("C108").PasteSpecial(-4104) | Out-Null
$ws.Range("D107").Copy() | Out-Null
$ws.Range("D108").PasteSpecial(-4104) | Out-Null

# Column E ("data_partida") holds the match's date/time serial with the
# custom date-time style (s="2"). Same format-copy trick as column A.
$ws.Range("E107").Copy() | Out-Null
$ws.Range("E108").PasteSpecial(-4122) | Out-Null
$ws.Range("E108").Value = 45233.91666666666

$ws.Range("F108").Value = "Guairena"
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = "Olimpia Asuncion"
$ws.Range("I108").Value = 1
$ws.Range("J108").Value = 4.77
$ws.Range("K108").Value = "29/10/2023 21:12"
$ws.Range("L108").Value = 5.26
$ws.Range("M108").Value = "03/11/2023 21:52"
$ws.Range("N108").Value = 3.9
$ws.Range("O108").Value = "29/10/2023 21:12"
$ws.Range("P108").Value = 3.81
$ws.Range("Q108").Value = "03/11/2023 21:55"
$ws.Range("R108").Value = 1.72
$ws.Range("S108").Value = "29/10/2023 21:12"
$ws.Range("T108").Value = 1.71
$ws.Range("U108").Value = "03/11/2023 21:54"
$ws.Range("V108").Value = "https://www.betexplorer.com/football/paraguay/primera-division/guairena-fc-olimpia-asuncion/pfDa4cTs/"

# --- Row 109 -------------------------------------------------------------

$ws.Range("A107").Copy() | Out-Null
$ws.Range("A109").PasteSpecial(-4122) | Out-Null
$ws.Range("A109").Value = 108

$ws.Range("B107").Copy() | Out-Null
$ws.Range("B109").PasteSpecial(-4104) | Out-Null
$ws.Range("C107").Copy() | Out-Null
$ws.Range("C109").PasteSpecial(-4104) | Out-Null
$ws.Range("D107").Copy() | Out-Null
$ws.Range("D109").PasteSpecial(-4104) | Out-Null

$ws.Range("E107").Copy() | Out-Null
$ws.Range("E109").PasteSpecial(-4122) | Out-Null
$ws.Range("E109").Value = 45234.02083333334

$ws.Range("F109").Value = "Cerro Porteno"
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = "Sp. Luqueno"
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 1.63
$ws.Range("K109").Value = "30/10/2023 22:12"
$ws.Range("L109").Value = 1.78
$ws.Range("M109").Value = "04/11/2023 00:26"
$ws.Range("N109").Value = 4.25
$ws.Range("O109").Value = "30/10/2023 22:12"
$ws.Range("P109").Value = 3.8
$ws.Range("Q109").Value = "04/11/2023 00:26"
$ws.Range("R109").Value = 4.6
$ws.Range("S109").Value = "30/10/2023 22:12"
$ws.Range("T109").Value = 4.77
$ws.Range("U109").Value = "04/11/2023 00:20"
$ws.Range("V109").Value = "https://www.betexplorer.com/football/paraguay/primera-division/cerro-porteno-sp-luqueno/SI243Hrl/"

$excel.CutCopyMode = $false
